$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 (VOLTAREN 75MG/3ML 3 AMP.) figures updated
$ws.Range("H24").Value = "1:1"
$ws.Range("P24").NumberFormat = "@"
$ws.Range("P24").Value = "33.6600"
$ws.Range("P24").NumberFormat = "0.00"
$ws.Range("Q24").Value = "0:2"

# Row 25 (سرنجات 3 سم) figures updated
$ws.Range("P25").NumberFormat = "@"
$ws.Range("P25").Value = "24.0000"
$ws.Range("P25").NumberFormat = "0.00"
$ws.Range("Q25").Value = "12:0"

# Recomputed total for column P
$ws.Range("P28").Value = 936.865

# Report timestamp refreshed
$ws.Range("A29").Value = "Wednesday, 6 August, 2025 11:22 AM"
